# Updated mean summary results
# The "total" rows (rows with no habitat_type) for each species are removed:
#  - the Largemouth Bass aggregate row (originally row 7)
#  - the Northern Pike aggregate row (originally row 13)
# Removing them shifts the remaining Northern Pike per-habitat rows up so
# they occupy rows 7-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the Largemouth Bass total row (row 7). Everything below shifts up
# by one, so the old row 13 (Northern Pike total row) is now row 12.
$ws.Rows.Item(7).Delete()

# Delete the Northern Pike total row, now at row 12.
$ws.Rows.Item(12).Delete()
